# Fix "not printing by month" issue:
#   - Mes (month) label changes from "July" to "June"
#   - The three summary totals (Visitantes/Acertos/Erros) are zeroed out
#   - The four "ranking" sub-tables (Respostas / Cidades / Notas / Idades)
#     lose all their data rows, leaving just the section title + header
#     row, and the sheet shrinks accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / summary block
$ws.Range("B2").Value = "June"
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0

# Row 6 stays blank (separator).
# Row 7: "Respostas Mais Acertadas" title, row 8: "Resposta" / "Total" header - unchanged.

# Delete the old data rows 9-33 entirely, then rebuild the trimmed layout
# below the "Resposta"/"Total" header (row 8).
$ws.Rows("9:33").Delete()

$ws.Range("A10").Value = "Cidades com Melhor Desempenho"
$ws.Range("A10").Font.Bold = $true

$ws.Range("A11").Value = "Cidade"
$ws.Range("B11").Value = "Total de Acertos"

$ws.Range("A13").Value = "Notas Mais Dadas"
$ws.Range("A13").Font.Bold = $true

$ws.Range("A14").Value = "Nota"
$ws.Range("B14").Value = "Total"

$ws.Range("A16").Value = "Idades Mais Visitadas"
$ws.Range("A16").Font.Bold = $true

$ws.Range("A17").Value = "Idade"
$ws.Range("B17").Value = "Total"
